$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-11-23"

# Update the row label for November to reflect the new "through" date
$ws.Range("A12").Value = "November (through 11-23)"

# Update the 2022 column (I) value for October row (row 11)
$ws.Range("I11").Value = 124

# Update the November row (row 12) with new data for 2015-2022
$ws.Range("B12").Value = 23
$ws.Range("C12").Value = 57
$ws.Range("D12").Value = 89
$ws.Range("E12").Value = 46
$ws.Range("F12").Value = 40
$ws.Range("G12").Value = 162
$ws.Range("H12").Value = 158
$ws.Range("I12").Value = 87

# Update the Total row (row 13) with new totals for 2015-2022
$ws.Range("B13").Value = 281
$ws.Range("C13").Value = 543
$ws.Range("D13").Value = 799
$ws.Range("E13").Value = 661
$ws.Range("F13").Value = 522
$ws.Range("G13").Value = 1219
$ws.Range("H13").Value = 1599
$ws.Range("I13").Value = 1485
